# Update countries & provincias Spain
# - Moldavia moves above Argelia (with updated Moldavia numbers)
# - Birmania moves above Guadalupe (with updated Birmania numbers)
# - San Martin (Parte Holandesa) moves above Guyana (with updated numbers)
# - Estados Unidos / Finlandia totals refreshed
# - "Datos actualizados" timestamp bumped to 16:22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 16:22"

# Estados Unidos (row 4) — refreshed totals
$ws.Range("B4").Value = 1012147
$ws.Range("C4").Value = 1791
$ws.Range("D4").Value = 139419
$ws.Range("E4").Value = 815795
$ws.Range("G4").Value = 136
$ws.Range("H4").Value = 56933

# Finlandia (row 54) — refreshed totals
$ws.Range("D54").Value = 2800
$ws.Range("E54").Value = 1741
$ws.Range("G54").Value = 6
$ws.Range("H54").Value = 199

# Moldavia now ranks above Argelia (row 58 becomes Moldavia with new data,
# row 59 becomes Argelia keeping its previous data)
$ws.Range("A58").Value = "Moldavia"
$ws.Range("B58").Value = 3638
$ws.Range("C58").Value = 157
$ws.Range("D58").Value = 925
$ws.Range("E58").Value = 2610
$ws.Range("F58").Value = 212
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 103

$ws.Range("A59").Value = "Argelia"
$ws.Range("B59").Value = 3517
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 1558
$ws.Range("E59").Value = 1527
$ws.Range("F59").Value = 40
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 432

# Birmania now ranks above Guadalupe (row 136 becomes Birmania with new
# data, row 137 becomes Guadalupe keeping its previous data)
$ws.Range("A136").Value = "Birmania"
$ws.Range("B136").Value = 149
$ws.Range("C136").Value = 3
$ws.Range("D136").Value = 16
$ws.Range("E136").Value = 128
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 5

$ws.Range("A137").Value = "Guadalupe"
$ws.Range("B137").Value = 149
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 82
$ws.Range("E137").Value = 55
$ws.Range("F137").Value = 11
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 12

# San Martin (Parte Holandesa) now ranks above Guyana (row 159 becomes San
# Martin with new data, row 160 becomes Guyana keeping its previous data)
$ws.Range("A159").Value = "San Martin (Parte Holandesa)"
$ws.Range("B159").Value = 75
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 33
$ws.Range("E159").Value = 29
$ws.Range("F159").Value = 7
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 13

$ws.Range("A160").Value = "Guyana"
$ws.Range("B160").Value = 74
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 15
$ws.Range("E160").Value = 51
$ws.Range("F160").Value = 5
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 8
